$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reset number format on B1 (will hold the new "Date and Time" text) back to the default style
$ws.Range("B1").ClearFormats()
# Apply the time-duration format (moved from old B1) to new B2
$ws.Range("B2").NumberFormat = "[hh]:mm:ss"

$ws.Cells.Item(1, 1).Value = "Date and Time"
$ws.Cells.Item(1, 2).Value = "2024-03-12 08:13:59.622000 to 2024-03-12 09:18:37.211000"

$ws.Cells.Item(2, 1).Value = "Total time taken for the ride"
$ws.Cells.Item(2, 2).Value = 0.04474824074074074

$ws.Cells.Item(3, 1).Value = "Actual Ampere-hours (Ah)"
$ws.Cells.Item(3, 2).Value = 31.87679666666666

$ws.Cells.Item(4, 1).Value = "Actual Watt-hours (Wh)"
$ws.Cells.Item(4, 2).Value = 1597.140168137222

$ws.Cells.Item(5, 1).Value = "Starting SoC (Ah)"
$ws.Cells.Item(5, 2).Value = 35.322

$ws.Cells.Item(6, 1).Value = "Ending SoC (Ah)"
$ws.Cells.Item(6, 2).Value = 3.45

$ws.Cells.Item(7, 1).Value = "Starting SoC (%)"
$ws.Cells.Item(7, 2).Value = 8

$ws.Cells.Item(8, 1).Value = "Ending SoC (%)"
$ws.Cells.Item(8, 2).Value = 89

$ws.Cells.Item(9, 1).Value = "Total distance covered (km)"
$ws.Cells.Item(9, 2).Value = 35.2506607867321

$ws.Cells.Item(10, 1).Value = "Total energy consumption(WH/KM)"
$ws.Cells.Item(10, 2).Value = 45.3080916071328

$ws.Cells.Item(11, 1).Value = "Total SOC consumed(%)"
$ws.Cells.Item(11, 2).Value = 81

$ws.Cells.Item(12, 1).Value = "Mode"
$ws.Cells.Item(12, 2).Value = "Custom mode`n52.47%`nEco mode`n42.89%"

$ws.Cells.Item(13, 1).Value = "Peak Power(kW)"
$ws.Cells.Item(13, 2).Value = 5429.287385000001

$ws.Cells.Item(14, 1).Value = "Average Power(kW)"
$ws.Cells.Item(14, 2).Value = -1497.318907628646

$ws.Cells.Item(15, 1).Value = "Total Energy Regenerated(kWh)"
$ws.Cells.Item(15, 2).Value = 0.60077226

$ws.Cells.Item(16, 1).Value = "Regenerative Effectiveness(%)"
$ws.Cells.Item(16, 2).Value = 0.0376013560653105

$ws.Cells.Item(17, 1).Value = "Highest Cell Voltage(V)"
$ws.Cells.Item(17, 2).Value = 3.371

$ws.Cells.Item(18, 1).Value = "Lowest Cell Voltage(V)"
$ws.Cells.Item(18, 2).Value = 2.98

$ws.Cells.Item(19, 1).Value = "Difference in Cell Voltage(V)"
$ws.Cells.Item(19, 2).Value = 0.391

$ws.Cells.Item(20, 1).Value = "Minimum Temperature(C)"
$ws.Cells.Item(20, 2).Value = 26

$ws.Cells.Item(21, 1).Value = "Maximum Temperature(C)"
$ws.Cells.Item(21, 2).Value = 39

$ws.Cells.Item(22, 1).Value = "Difference in Temperature(C)"
$ws.Cells.Item(22, 2).Value = 13

$ws.Cells.Item(23, 1).Value = "Maximum Fet Temperature-BMS(C)"
$ws.Cells.Item(23, 2).Value = 63

$ws.Cells.Item(24, 1).Value = "Maximum Afe Temperature-BMS(C)"
$ws.Cells.Item(24, 2).Value = 58

$ws.Cells.Item(25, 1).Value = "Maximum PCB Temperature-BMS(C)"
$ws.Cells.Item(25, 2).Value = 59

$ws.Cells.Item(26, 1).Value = "Maximum MCU Temperature(C)"
$ws.Cells.Item(26, 2).Value = 54

$ws.Cells.Item(27, 1).Value = "Maximum Motor Temperature(C)"
$ws.Cells.Item(27, 2).Value = 89

$ws.Cells.Item(28, 1).Value = "Abnormal Motor Temperature Detected(C)"
$ws.Cells.Item(28, 2).Value = 0

$ws.Cells.Item(29, 1).Value = "highest cell temp(C)"
$ws.Cells.Item(29, 2).Value = 39

$ws.Cells.Item(30, 1).Value = "lowest cell temp(C)"
$ws.Cells.Item(30, 2).Value = 26

$ws.Cells.Item(31, 1).Value = "Difference between Highest and Lowest Cell Temperature at 100% SOC(C)"
$ws.Cells.Item(31, 2).Value = 13

$ws.Cells.Item(32, 1).Value = "Battery Voltage(V)"
$ws.Cells.Item(32, 2).Value = 53

$ws.Cells.Item(33, 1).Value = "Total energy charged(kWh)"
$ws.Cells.Item(33, 2).Value = 1.689470223333333

$ws.Cells.Item(34, 1).Value = "Electricity consumption units(kW)"
$ws.Cells.Item(34, 2).Value = [double]"1.213909167768389e-07"

$ws.Cells.Item(35, 1).Value = "Cycle Count of battery"
$ws.Cells.Item(35, 2).Value = 39

$ws.Cells.Item(36, 1).Value = "Idling time percentage"
$ws.Cells.Item(36, 2).Value = 22.70450751252087

$ws.Cells.Item(37, 1).Value = "Time spent in 0-10 km/h"
$ws.Cells.Item(37, 2).Value = 2.899276572064552

$ws.Cells.Item(38, 1).Value = "Time spent in 10-20 km/h"
$ws.Cells.Item(38, 2).Value = 4.774624373956595

$ws.Cells.Item(39, 1).Value = "Time spent in 20-30 km/h"
$ws.Cells.Item(39, 2).Value = 9.315525876460768

$ws.Cells.Item(40, 1).Value = "Time spent in 30-40 km/h"
$ws.Cells.Item(40, 2).Value = 30.28380634390651

$ws.Cells.Item(41, 1).Value = "Time spent in 40-50 km/h"
$ws.Cells.Item(41, 2).Value = 5.779076238174736

$ws.Cells.Item(42, 1).Value = "Time spent in 50-60 km/h"
$ws.Cells.Item(42, 2).Value = 8.597662771285476

$ws.Cells.Item(43, 1).Value = "Time spent in 60-70 km/h"
$ws.Cells.Item(43, 2).Value = 11.47746243739566

$ws.Cells.Item(44, 1).Value = "Time spent in 70-80 km/h"
$ws.Cells.Item(44, 2).Value = 4.067890929326656

$ws.Cells.Item(45, 1).Value = "Time spent in 80-90 km/h"
$ws.Cells.Item(45, 2).Value = 0

# Re-fit row height for the multi-line "Mode" cell so no stray explicit height is stored
$ws.Rows.Item(12).AutoFit()
